# Updates cryptos list figures (Price / Volume(1h)) to the latest snapshot.
# Rows 44/45 (OKB <-> Bittensor) also swap places in this snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.262.87"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "3.793.23"
$ws.Range("E3").Value = "  +1.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.80"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.04"
$ws.Range("E6").Value = "  -2.73%  "

# Row 7
$ws.Range("D7").Value = "3.787.86"
$ws.Range("E7").Value = "  +1.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  +3.10%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.33"
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -0.45%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.34"
$ws.Range("E13").Value = "  -2.65%  "

# Row 14
$ws.Range("E14").Value = "  -0.69%  "

# Row 15
$ws.Range("D15").Value = "4.428.89"
$ws.Range("E15").Value = "  +1.18%  "

# Row 16
$ws.Range("D16").Value = "3.798.63"
$ws.Range("E16").Value = "  +1.51%  "

# Row 17
$ws.Range("D17").Value = "69.345.19"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.45"
$ws.Range("E18").Value = "  +2.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.47"
$ws.Range("E19").Value = "  +2.27%  "

# Row 20
$ws.Range("E20").Value = "  -0.47%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  +4.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.26"
$ws.Range("E22").Value = "  -0.46%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("E23").Value = "  -0.53%  "

# Row 24
$ws.Range("E24").Value = "  -2.14%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.83"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  -2.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -2.41%  "

# Row 29
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("E30").Value = "  +0.05%  "

# Row 31
$ws.Range("E31").Value = "  +2.11%  "

# Row 32
$ws.Range("E32").Value = "  -5.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.05"
$ws.Range("E33").Value = "  +0.52%  "

# Row 34
$ws.Range("D34").Value = "3.943.59"
$ws.Range("E34").Value = "  +1.21%  "

# Row 35
$ws.Range("D35").Value = "3.741.19"
$ws.Range("E35").Value = "  +1.54%  "

# Row 36
$ws.Range("E36").Value = "  -1.30%  "

# Row 37
$ws.Range("E37").Value = "  +6.19%  "

# Row 38
$ws.Range("E38").Value = "  +0.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.94"
$ws.Range("E39").Value = "  +1.17%  "

# Row 40
$ws.Range("E40").Value = "  +0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  +2.42%  "

# Row 43
$ws.Range("E43").Value = "  +0.85%  "

# Row 44
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "425.98"
$ws.Range("E44").Value = "  -3.02%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.44"
$ws.Range("E45").Value = "  -1.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.42"
$ws.Range("E46").Value = "  -0.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.34"
$ws.Range("E48").Value = "  +0.52%  "

# Row 49
$ws.Range("D49").Value = "2.816.51"
$ws.Range("E49").Value = "  +1.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.91"
$ws.Range("E50").Value = "  -1.60%  "

# Row 51
$ws.Range("E51").Value = "  +5.51%  "
